$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Underline the whole first paragraph ("MSP 430 To do list"),
#    including the paragraph mark (so the pPr/rPr also gets w:u).
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Font.Underline = 1

# ------------------------------------------------------------------
# 2) Remove the old (hidden) "_GoBack" bookmark that currently sits
#    near "-Fix the 3.3 label to 5V coming off of the USB" so the
#    name is free to be reused at the top of the document.
# ------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
if ($oldGoBack -ne $null) {
    $oldGoBack.Delete()
}

# ------------------------------------------------------------------
# 3) Delete the three paragraphs that were removed from the "to do"
#    list: "Incorporate ferrite bead footprint", "Resitor footprint"
#    and "Capacitor Footprint".
# ------------------------------------------------------------------
$delStart = $d.Paragraphs.Item(2).Range.Start
$delEnd = $d.Paragraphs.Item(4).Range.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

# ------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark spanning the first paragraph
#    (including its paragraph mark) now that the list has been
#    trimmed down, matching the document's final, saved cursor spot.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$d.Bookmarks.Add("_GoBack", $p1.Range)
